$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data one column to the right to make room for a new
# index column (mirrors a pandas reset_index() export).
$ws.Columns.Item(1).Insert()

# New index column A: numeric row index (0-based), same header style as
# the other header cells (bold/centered/bordered).
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5

$ws.Range("B1").Copy()
$ws.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "S&P 500" summary row values changed.
$ws.Range("C7").Value = 24.12
$ws.Range("D7").Value = 0.13
$ws.Range("E7").Value = 1.55
